$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + "27.558.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'" + "  +0.64%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'" + "1.641.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'" + "  -0.73%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'" + "  -0.03%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'" + "212.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'" + "  -0.43%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'" + "0.536"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'" + "  +4.57%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'" + "  -0.04%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'" + "22.92"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'" + "  -4.19%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'" + "  -1.53%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "'" + "  -0.52%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'" + "0.0889"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'" + "  +1.33%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'" + "1.874.87"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'" + "  -0.77%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'" + "1.647.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'" + "  -0.30%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("E14").Value = "'" + "  -0.89%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("E15").Value = "'" + "  -1.38%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("E16").Value = "'" + "  -2.44%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'" + "27.551.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'" + "  +0.53%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'" + "228.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'" + "  -1.27%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'" + "0.0₃0723"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'" + "  -0.22%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'" + "7.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'" + "  +2.11%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "'" + "  -0.02%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'" + "  -1.63%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'" + "10.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'" + "  +7.61%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'" + "  -3.64%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'" + "149.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'" + "  +1.71%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'" + "6.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'" + "  -3.12%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'" + "  +1.42%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'" + "  -0.09%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'" + "  -1.62%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "'" + "  -0.95%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'" + "0.0485"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'" + "  -2.37%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = "'" + "  -0.03%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'" + "3.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'" + "  +1.74%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'" + "1.427.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'" + "  -2.27%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'" + "1.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'" + "  +2.23%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D37").Value = "'" + "0.573"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'" + "  +0.41%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'" + "  -3.61%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = "'" + "  -1.29%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'" + "0.890"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'" + "  +13.69%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'" + "  -2.22%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'" + "  -0.01%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("B43").Value = "'" + "MXToken"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'" + "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'" + "2.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'" + "  +2.13%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("B44").Value = "'" + "mCoin"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'" + "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'" + "2.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'" + "  -1.17%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("B45").Value = "'" + "FraxShare"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'" + "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'" + "5.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'" + "  +1.07%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'" + "65.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'" + "  -0.26%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'" + "1.783.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'" + "  -0.74%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'" + "1.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'" + "  -2.53%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'" + "  -2.17%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'" + "  +1.00%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'" + "0.0984"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'" + "  -2.56%  "
$ws.Range("E51").Style = "Normal"
